# Update the monthly salary input ("Salário") on the simulator sheet.
# Typing a new value into D4 and pressing Enter would normally move the
# active cell down to D5 - we reproduce that final selection state too.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")
$ws.Activate()

$ws.Range("D4").Value = 5000

$ws.Range("D5").Select()
